# Auto-generated edit script: updates crypto price/volume table
# to match the GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.656.81"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "2.289.81"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "115.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E7").Value = "  +2.61%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0943"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").Value = "2.633.01"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.885"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("D17").Value = "2.293.74"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "43.696.00"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0914"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("E36").Value = "  +5.42%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.108"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("E40").Value = "  +7.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("E45").Value = "  -7.16%  "
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("E47").Value = "  +4.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +33.16%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.14%  "
